$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A59").Value = "Davide Scarperi"
$ws.Range("B59").Value = "Alessandro  Ruele  | FC GORILLAZ"
$ws.Range("C59").Value = "Andrea Conzatti | FC SAVIGNANO"
$ws.Range("D59").Value = "Alessio Bragagna | SHARK ATTACK"
$ws.Range("E59").Value = "Daniel Pedrotti | IMONTAGNA"
$ws.Range("F59").Value = "Emanuele  valduga | wanda tim"
